$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Bump the CodeSystem version
$ws.Range("B3").Value = "6.0.0"

# 2. Update the publication date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# 3. Publisher now has a real value
$ws.Range("B9").Value = "Alvearie Team"

# 4. Remove the old duplicated "Contact" / "No display for ContactDetail" rows
#    (rows 10 and 11 both held the same Contact info before the edit).
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()

# 5. Insert a single replacement row for the new "Jurisdiction" property.
#    Insert the blank row first, then copy the formatting of the row above
#    into it, so the resulting cell style matches the rest of the table.
$ws.Rows.Item(10).Insert()
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# 6. "Case Sensitive" now has an explicit value of true (it shifted up one
#    row to row 14 after the net row-count change above). Entering the bare
#    word "true" would be auto-coerced to a Boolean by Excel, but the source
#    workbook stores it as plain text, so round-trip it through a nearby
#    helper cell (as a text-formula result) and paste just the value back.
$ws.Range("C14").Formula = "=""true"""
$ws.Range("C14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("C14").ClearContents()

Write-Output "edits applied"
